$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 71
$ws.Range("F3").Value = 1784
$ws.Range("F6").Value = 1160
$ws.Range("F10").Value = 22
$ws.Range("F11").Value = 1532
$ws.Range("F12").Value = 3144
$ws.Range("F13").Value = 689
$ws.Range("F14").Value = 1851
$ws.Range("F15").Value = 1835
$ws.Range("F16").Value = 893
$ws.Range("F17").Value = 305
$ws.Range("F19").Value = 1517
$ws.Range("F20").Value = 309
$ws.Range("F23").Value = 1306
$ws.Range("F24").Value = 426
$ws.Range("F25").Value = 521
$ws.Range("F26").Value = 202
$ws.Range("F27").Value = 7755
$ws.Range("F28").Value = 7177
$ws.Range("F29").Value = 779
$ws.Range("F30").Value = 597
$ws.Range("F31").Value = 1721
$ws.Range("F32").Value = 98
$ws.Range("F33").Value = 250

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 27

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 52

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 71
$ws.Range("F3").Value = 52
$ws.Range("F5").Value = 1784
$ws.Range("F8").Value = 1160
$ws.Range("F13").Value = 22
$ws.Range("F14").Value = 1532
$ws.Range("F15").Value = 3144
$ws.Range("F16").Value = 689
$ws.Range("F17").Value = 1851
$ws.Range("F18").Value = 1835
$ws.Range("F19").Value = 893
$ws.Range("F20").Value = 305
$ws.Range("F22").Value = 1517
$ws.Range("F23").Value = 309
$ws.Range("F28").Value = 1306
$ws.Range("F29").Value = 426
$ws.Range("F30").Value = 521
$ws.Range("F31").Value = 202
$ws.Range("F32").Value = 7755
$ws.Range("F33").Value = 7177
$ws.Range("F34").Value = 779
$ws.Range("F35").Value = 597
$ws.Range("F36").Value = 1721
$ws.Range("F38").Value = 27
$ws.Range("F39").Value = 98
$ws.Range("F40").Value = 250
